# Applies NATMI recompute update: Ligand/Receptor-expressing cells count changed from 1 to 3,
# which propagates through total-expression, specificity, and edge-weight columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 162.399297
$ws.Range("H2").Value = 487.197891
$ws.Range("I2").Value = 0.3910371682630009
$ws.Range("J2").Value = 0.3910371682630009
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 490.031855
$ws.Range("N2").Value = 1470.095565
$ws.Range("O2").Value = 0.6686419015677429
$ws.Range("P2").Value = 0.6686419015677431
$ws.Range("Q2").Value = 79580.82875960594
$ws.Range("R2").Value = 716227.4588364534
$ws.Range("S2").Value = 0.2614638357710384
$ws.Range("T2").Value = 0.2614638357710384

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 162.399297
$ws.Range("H3").Value = 487.197891
$ws.Range("I3").Value = 0.3910371682630009
$ws.Range("J3").Value = 0.3910371682630009
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 62.79306433333334
$ws.Range("N3").Value = 188.379193
$ws.Range("O3").Value = 0.0856802950924601
$ws.Range("P3").Value = 0.08568029509246011
$ws.Range("Q3").Value = 10197.54950420911
$ws.Range("R3").Value = 91777.94553788197
$ws.Range("S3").Value = 0.03350417996889389
$ws.Range("T3").Value = 0.0335041799688939

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 162.399297
$ws.Range("H4").Value = 487.197891
$ws.Range("I4").Value = 0.3910371682630009
$ws.Range("J4").Value = 0.3910371682630009
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.5977846666666666
$ws.Range("N4").Value = 1.793354
$ws.Range("O4").Value = 0.0008156691696053909
$ws.Range("P4").Value = 0.000815669169605391
$ws.Range("Q4").Value = 97.07980962404599
$ws.Range("R4").Value = 873.7182866164139
$ws.Range("S4").Value = 0.0003189569623219255
$ws.Range("T4").Value = 0.0003189569623219256

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 162.399297
$ws.Range("H5").Value = 487.197891
$ws.Range("I5").Value = 0.3910371682630009
$ws.Range("J5").Value = 0.3910371682630009
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 179.453674
$ws.Range("N5").Value = 538.361022
$ws.Range("O5").Value = 0.2448621341701915
$ws.Range("P5").Value = 0.2448621341701916
$ws.Range("Q5").Value = 29143.15050166718
$ws.Range("R5").Value = 262288.3545150046
$ws.Range("S5").Value = 0.0957501955607467
$ws.Range("T5").Value = 0.09575019556074671

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 65.41736466666667
$ws.Range("H6").Value = 196.252094
$ws.Range("I6").Value = 0.1575168212364948
$ws.Range("J6").Value = 0.1575168212364948
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 490.031855
$ws.Range("N6").Value = 1470.095565
$ws.Range("O6").Value = 0.6686419015677429
$ws.Range("P6").Value = 0.6686419015677431
$ws.Range("Q6").Value = 32056.59255681813
$ws.Range("R6").Value = 288509.3330113631
$ws.Range("S6").Value = 0.1053223468804761
$ws.Range("T6").Value = 0.1053223468804761

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 65.41736466666667
$ws.Range("H7").Value = 196.252094
$ws.Range("I7").Value = 0.1575168212364948
$ws.Range("J7").Value = 0.1575168212364948
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 62.79306433333334
$ws.Range("N7").Value = 188.379193
$ws.Range("O7").Value = 0.0856802950924601
$ws.Range("P7").Value = 0.08568029509246011
$ws.Range("Q7").Value = 4107.756788031128
$ws.Range("R7").Value = 36969.81109228014
$ws.Range("S7").Value = 0.01349608772556916
$ws.Range("T7").Value = 0.01349608772556916

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 65.41736466666667
$ws.Range("H8").Value = 196.252094
$ws.Range("I8").Value = 0.1575168212364948
$ws.Range("J8").Value = 0.1575168212364948
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.5977846666666666
$ws.Range("N8").Value = 1.793354
$ws.Range("O8").Value = 0.0008156691696053909
$ws.Range("P8").Value = 0.000815669169605391
$ws.Range("Q8").Value = 39.10549753147511
$ws.Range("R8").Value = 351.949477783276
$ws.Range("S8").Value = 0.0001284816147768525
$ws.Range("T8").Value = 0.0001284816147768526

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 65.41736466666667
$ws.Range("H9").Value = 196.252094
$ws.Range("I9").Value = 0.1575168212364948
$ws.Range("J9").Value = 0.1575168212364948
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 179.453674
$ws.Range("N9").Value = 538.361022
$ws.Range("O9").Value = 0.2448621341701915
$ws.Range("P9").Value = 0.2448621341701916
$ws.Range("Q9").Value = 11739.38643283112
$ws.Range("R9").Value = 105654.4778954801
$ws.Range("S9").Value = 0.03856990501567267
$ws.Range("T9").Value = 0.03856990501567267

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 126.3069433333333
$ws.Range("H10").Value = 378.92083
$ws.Range("I10").Value = 0.3041313008456065
$ws.Range("J10").Value = 0.3041313008456065
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 490.031855
$ws.Range("N10").Value = 1470.095565
$ws.Range("O10").Value = 0.6686419015677429
$ws.Range("P10").Value = 0.6686419015677431
$ws.Range("Q10").Value = 61894.42574101322
$ws.Range("R10").Value = 557049.8316691191
$ws.Range("S10").Value = 0.2033549313236776
$ws.Range("T10").Value = 0.2033549313236777

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 126.3069433333333
$ws.Range("H11").Value = 378.92083
$ws.Range("I11").Value = 0.3041313008456065
$ws.Range("J11").Value = 0.3041313008456065
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 62.79306433333334
$ws.Range("N11").Value = 188.379193
$ws.Range("O11").Value = 0.0856802950924601
$ws.Range("P11").Value = 0.08568029509246011
$ws.Range("Q11").Value = 7931.200018476688
$ws.Range("R11").Value = 71380.8001662902
$ws.Range("S11").Value = 0.02605805960330532
$ws.Range("T11").Value = 0.02605805960330533

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 126.3069433333333
$ws.Range("H12").Value = 378.92083
$ws.Range("I12").Value = 0.3041313008456065
$ws.Range("J12").Value = 0.3041313008456065
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.5977846666666666
$ws.Range("N12").Value = 1.793354
$ws.Range("O12").Value = 0.0008156691696053909
$ws.Range("P12").Value = 0.000815669169605391
$ws.Range("Q12").Value = 75.50435401820222
$ws.Range("R12").Value = 679.53918616382
$ws.Range("S12").Value = 0.0002480705256117432
$ws.Range("T12").Value = 0.0002480705256117432

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 126.3069433333333
$ws.Range("H13").Value = 378.92083
$ws.Range("I13").Value = 0.3041313008456065
$ws.Range("J13").Value = 0.3041313008456065
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 179.453674
$ws.Range("N13").Value = 538.361022
$ws.Range("O13").Value = 0.2448621341701915
$ws.Range("P13").Value = 0.2448621341701916
$ws.Range("Q13").Value = 22666.24503287648
$ws.Range("R13").Value = 203996.2052958883
$ws.Range("S13").Value = 0.07447023939301177
$ws.Range("T13").Value = 0.0744702393930118

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 61.180387
$ws.Range("H14").Value = 183.541161
$ws.Range("I14").Value = 0.1473147096548978
$ws.Range("J14").Value = 0.1473147096548978
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 490.031855
$ws.Range("N14").Value = 1470.095565
$ws.Range("O14").Value = 0.6686419015677429
$ws.Range("P14").Value = 0.6686419015677431
$ws.Range("Q14").Value = 29980.33853122788
$ws.Range("R14").Value = 269823.046781051
$ws.Range("S14").Value = 0.09850078759255081
$ws.Range("T14").Value = 0.09850078759255083

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 61.180387
$ws.Range("H15").Value = 183.541161
$ws.Range("I15").Value = 0.1473147096548978
$ws.Range("J15").Value = 0.1473147096548978
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 62.79306433333334
$ws.Range("N15").Value = 188.379193
$ws.Range("O15").Value = 0.0856802950924601
$ws.Range("P15").Value = 0.08568029509246011
$ws.Range("Q15").Value = 3841.703976829231
$ws.Range("R15").Value = 34575.33579146308
$ws.Range("S15").Value = 0.01262196779469172
$ws.Range("T15").Value = 0.01262196779469173

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 61.180387
$ws.Range("H16").Value = 183.541161
$ws.Range("I16").Value = 0.1473147096548978
$ws.Range("J16").Value = 0.1473147096548978
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.5977846666666666
$ws.Range("N16").Value = 1.793354
$ws.Range("O16").Value = 0.0008156691696053909
$ws.Range("P16").Value = 0.000815669169605391
$ws.Range("Q16").Value = 36.57269724933266
$ws.Range("R16").Value = 329.1542752439939
$ws.Range("S16").Value = 0.0001201600668948697
$ws.Range("T16").Value = 0.0001201600668948698

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 61.180387
$ws.Range("H17").Value = 183.541161
$ws.Range("I17").Value = 0.1473147096548978
$ws.Range("J17").Value = 0.1473147096548978
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 179.453674
$ws.Range("N17").Value = 538.361022
$ws.Range("O17").Value = 0.2448621341701915
$ws.Range("P17").Value = 0.2448621341701916
$ws.Range("Q17").Value = 10979.04522389184
$ws.Range("R17").Value = 98811.40701502655
$ws.Range("S17").Value = 0.03607179420076039
$ws.Range("T17").Value = 0.03607179420076041

